$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The TreatmentTab query (row 5, column B) had a redundant CONCAT(...) wrapper
# around the REPLACE(...) call used to build the "Treatment Agent" column.
# Fix it by removing the unnecessary CONCAT() wrapper, keeping the rest of the
# query (and the cell's existing wrap/size formatting) untouched.
$oldQuery = $ws.Range("B5").Value()
$newQuery = $oldQuery.Replace(
    'CONCAT(REPLACE(trt.treatment_agent, '';'', '', '')) AS "Treatment Agent",',
    'REPLACE(trt.treatment_agent, '';'', '', '') AS "Treatment Agent",'
)
$ws.Range("B5").Value = $newQuery

# Reflect the author's saved view: scrolled/selected on C5 instead of B7.
$ws.Range("C5").Select()
